$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.568.84"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.790.05"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.39"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.554"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.89%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.27"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.280"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0687"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.046.62"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.781.34"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.89"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.585.12"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.629"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.29"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.38"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.17"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0788"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +6.79%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.15"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.90"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.37"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.43%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.06"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.114"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0516"
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.74"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.430.19"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.01%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.633"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0189"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "83.54"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.73%  "
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.903"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.06"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.95"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.37%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0496"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -3.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.945.01"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.90%  "
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.97"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.34%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.38"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.34%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "49.54"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.57%  "
